$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E6").Value = 31
$ws.Range("E7").Value = 25
$ws.Range("E10").Value = 298
$ws.Range("F10").Value = 146
$ws.Range("H10").Value = 146
$ws.Range("E11").Value = 209
$ws.Range("F11").Value = 124
$ws.Range("H11").Value = 124
$ws.Range("E12").Value = 306
$ws.Range("F12").Value = 179
$ws.Range("H12").Value = 179
$ws.Range("F13").Value = 47
$ws.Range("H13").Value = 47
$ws.Range("E14").Value = 85
$ws.Range("E16").Value = 127
$ws.Range("F16").Value = 73
$ws.Range("H16").Value = 73
$ws.Range("E17").Value = 57
$ws.Range("E18").Value = 42
$ws.Range("E20").Value = 65
$ws.Range("E21").Value = 91
$ws.Range("F21").Value = 54
$ws.Range("H21").Value = 54
$ws.Range("E23").Value = 119
$ws.Range("F23").Value = 64
$ws.Range("H23").Value = 64
$ws.Range("E24").Value = 137
$ws.Range("F24").Value = 72
$ws.Range("H24").Value = 72
$ws.Range("E25").Value = 147
$ws.Range("F25").Value = 74
$ws.Range("H25").Value = 74
$ws.Range("E26").Value = 82
$ws.Range("F26").Value = 54
$ws.Range("H26").Value = 54
$ws.Range("E27").Value = 196
$ws.Range("E28").Value = 120
$ws.Range("F28").Value = 45
$ws.Range("H28").Value = 45
$ws.Range("E29").Value = 120
$ws.Range("F30").Value = 80
$ws.Range("H30").Value = 80
$ws.Range("E32").Value = 126
$ws.Range("F32").Value = 70
$ws.Range("H32").Value = 70
$ws.Range("E33").Value = 188
$ws.Range("F33").Value = 99
$ws.Range("H33").Value = 99
$ws.Range("E34").Value = 140
$ws.Range("F34").Value = 89
$ws.Range("H34").Value = 89
$ws.Range("E36").Value = 43
$ws.Range("F36").Value = 26
$ws.Range("H36").Value = 26
$ws.Range("E37").Value = 102
$ws.Range("E38").Value = 61
$ws.Range("F38").Value = 48
$ws.Range("H38").Value = 48
$ws.Range("E39").Value = 125
$ws.Range("E40").Value = 176
$ws.Range("E41").Value = 246
$ws.Range("E42").Value = 229
$ws.Range("F42").Value = 121
$ws.Range("H42").Value = 121
$ws.Range("E43").Value = 69
$ws.Range("F43").Value = 36
$ws.Range("H43").Value = 36
$ws.Range("E44").Value = 194
$ws.Range("F44").Value = 106
$ws.Range("H44").Value = 106
$ws.Range("E45").Value = 74
$ws.Range("F45").Value = 43
$ws.Range("H45").Value = 43
$ws.Range("E46").Value = 189
$ws.Range("F46").Value = 111
$ws.Range("H46").Value = 111
$ws.Range("E47").Value = 288
$ws.Range("F47").Value = 148
$ws.Range("H47").Value = 148
$ws.Range("E48").Value = 132
$ws.Range("E49").Value = 165
$ws.Range("F49").Value = 81
$ws.Range("H49").Value = 81
$ws.Range("E50").Value = 141
$ws.Range("E51").Value = 135
$ws.Range("F51").Value = 62
$ws.Range("H51").Value = 62
